$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '61.650.24'
Set-TextValue $ws.Range("E2") '  +0.86%  '

Set-TextValue $ws.Range("D3") '3.451.13'
Set-TextValue $ws.Range("E3") '  +1.56%  '

Set-TextValue $ws.Range("E4") '  -0.13%  '

Set-TextValue $ws.Range("D5") '577.52'
Set-TextValue $ws.Range("E5") '  +0.79%  '

Set-TextValue $ws.Range("D6") '145.78'
Set-TextValue $ws.Range("E6") '  +5.00%  '

Set-TextValue $ws.Range("D7") '3.451.65'
Set-TextValue $ws.Range("E7") '  +1.55%  '

Set-TextValue $ws.Range("E8") '  -0.02%  '

Set-TextValue $ws.Range("E9") '  +1.91%  '

Set-TextValue $ws.Range("E10") '  +0.09%  '

Set-TextValue $ws.Range("E11") '  +3.73%  '

Set-TextValue $ws.Range("E12") '  +2.84%  '

Set-TextValue $ws.Range("D13") '4.038.74'
Set-TextValue $ws.Range("E13") '  +1.43%  '

Set-TextValue $ws.Range("D14") '28.49'
Set-TextValue $ws.Range("E14") '  +7.23%  '

Set-TextValue $ws.Range("E15") '  -0.40%  '

Set-TextValue $ws.Range("D16") '0.0000174'
Set-TextValue $ws.Range("E16") '  +1.30%  '

Set-TextValue $ws.Range("D17") '3.446.78'
Set-TextValue $ws.Range("E17") '  +1.37%  '

Set-TextValue $ws.Range("D18") '61.762.55'
Set-TextValue $ws.Range("E18") '  +0.88%  '

Set-TextValue $ws.Range("D19") '6.39'
Set-TextValue $ws.Range("E19") '  +7.44%  '

Set-TextValue $ws.Range("D20") '14.35'
Set-TextValue $ws.Range("E20") '  +3.55%  '

Set-TextValue $ws.Range("D21") '9.43'
Set-TextValue $ws.Range("E21") '  +1.31%  '

Set-TextValue $ws.Range("D22") '403.22'
Set-TextValue $ws.Range("E22") '  +7.23%  '

Set-TextValue $ws.Range("E23") '  +2.92%  '

Set-TextValue $ws.Range("D24") '74.57'
Set-TextValue $ws.Range("E24") '  +4.80%  '

Set-TextValue $ws.Range("E25") '  +0.51%  '

Set-TextValue $ws.Range("E26") '  -0.56%  '

Set-TextValue $ws.Range("D27") '0.0000124'
Set-TextValue $ws.Range("E27") '  +1.50%  '

Set-TextValue $ws.Range("D28") '3.588.29'
Set-TextValue $ws.Range("E28") '  +1.61%  '

Set-TextValue $ws.Range("D29") '0.183'
Set-TextValue $ws.Range("E29") '  +4.42%  '

Set-TextValue $ws.Range("D30") '7.63'
Set-TextValue $ws.Range("E30") '  +2.85%  '

Set-TextValue $ws.Range("E31") '  +0.19%  '

Set-TextValue $ws.Range("D32") '8.26'
Set-TextValue $ws.Range("E32") '  +1.64%  '

Set-TextValue $ws.Range("E33") '  +2.10%  '

Set-TextValue $ws.Range("E34") '  -9.76%  '

Set-TextValue $ws.Range("E35") '  -0.09%  '

Set-TextValue $ws.Range("D36") '23.96'
Set-TextValue $ws.Range("E36") '  +2.23%  '

Set-TextValue $ws.Range("D37") '7.06'
Set-TextValue $ws.Range("E37") '  +2.68%  '

Set-TextValue $ws.Range("D38") '3.476.41'
Set-TextValue $ws.Range("E38") '  +1.63%  '

Set-TextValue $ws.Range("D39") '1.57'
Set-TextValue $ws.Range("E39") '  +0.01%  '

Set-TextValue $ws.Range("D40") '5.14'
Set-TextValue $ws.Range("E40") '  +0.55%  '

Set-TextValue $ws.Range("D41") '167.15'
Set-TextValue $ws.Range("E41") '  +0.61%  '

Set-TextValue $ws.Range("D42") '0.0792'
Set-TextValue $ws.Range("E42") '  +2.71%  '

Set-TextValue $ws.Range("D43") '27.24'
Set-TextValue $ws.Range("E43") '  +5.09%  '

Set-TextValue $ws.Range("D44") '0.803'
Set-TextValue $ws.Range("E44") '  +3.33%  '

Set-TextValue $ws.Range("B45") 'Stacks'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D45") '1.75'
Set-TextValue $ws.Range("E45") '  +0.35%  '

Set-TextValue $ws.Range("B46") 'Filecoin'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D46") '4.53'
Set-TextValue $ws.Range("E46") '  +2.87%  '

Set-TextValue $ws.Range("B47") 'OKB'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D47") '42.47'
Set-TextValue $ws.Range("E47") '  +1.26%  '

Set-TextValue $ws.Range("B48") 'FirstDigitalUSD'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D48") '0.999'
Set-TextValue $ws.Range("E48") '  -0.22%  '

Set-TextValue $ws.Range("D49") '2.609.06'
Set-TextValue $ws.Range("E49") '  +3.68%  '

Set-TextValue $ws.Range("E50") '  -1.53%  '

Set-TextValue $ws.Range("D51") '6.96'
Set-TextValue $ws.Range("E51") '  +2.71%  '
